$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial for every data row (2..265).
# All of them move forward by one day: 45188 -> 45189.
for ($r = 2; $r -le 265; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $val = $cell.Value2()
    if ($val -eq 45188) {
        $cell.Value = 45189
    }
}

# Rows 262, 263 and 265 have their "Beteckning" (col A) and "Area (ha)"
# (col G) values rotated: 265 -> 262 -> 263 -> 265.
$a262 = $ws.Cells.Item(262, 1).Value()
$a263 = $ws.Cells.Item(263, 1).Value()
$a265 = $ws.Cells.Item(265, 1).Value()

$g262 = $ws.Cells.Item(262, 7).Value2()
$g263 = $ws.Cells.Item(263, 7).Value2()
$g265 = $ws.Cells.Item(265, 7).Value2()

$ws.Cells.Item(262, 1).Value = $a265
$ws.Cells.Item(263, 1).Value = $a262
$ws.Cells.Item(265, 1).Value = $a263

$ws.Cells.Item(262, 7).Value = $g265
$ws.Cells.Item(263, 7).Value = $g262
$ws.Cells.Item(265, 7).Value = $g263
